# SDSS4-Slides-Part2 edit script
#
# Summary of the change being applied:
#  1. Duplicate slide 14 ("Open Discussion") and move the duplicate to sit
#     BEFORE the original (i.e. it becomes the new slide 14), then retitle
#     the duplicate's headline text box to "Questions about Infrastructure?".
#     This pushes the original "Open Discussion" slide to position 15, the
#     "Discussion prompts" slide to position 16 and the "Thank You" slide to
#     position 17.
#  2. Rewrite the three numbered discussion-prompt paragraphs on the
#     (now) 16th slide: reorder/merge the old four questions into three,
#     with several words emphasised (bold+italic).

$p = $ppt.ActivePresentation

# --- 1. Duplicate "Open Discussion" (slide 14) and reorder -----------------

$openDiscussion = $p.Slides.Item(14)
$dupRange = $openDiscussion.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(14)

# Retitle the duplicate (now slide 14) headline text box.
$titleShape = $p.Slides.Item(14).Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Questions about Infrastructure?"

# --- 2. Rewrite the discussion-prompt questions (now slide 16) -------------

$promptsSlide = $p.Slides.Item(16)
$promptsShape = $promptsSlide.Shapes.Item(3)
$tr = $promptsShape.TextFrame.TextRange

$tr.Text = "Can you share any successes, advice, or best practices introducing reproducibility and replicability in your spatial data science scholarship (research or teaching)?`rWhat barriers do you perceive to adopting open and reproducible research practices in your own scholarship?`rCould any resources, changes, or incentives help overcome those barriers?"

# Bold + italic emphasis ranges (1-indexed, absolute within the whole
# TextRange; paragraph marks count as one character each).
$emphasisRanges = @(
    @(19, 9),    # "successes"
    @(30, 6),    # "advice"
    @(41, 15),   # "best practices "
    @(172, 8),   # "barriers"
    @(285, 10),  # "resources,"
    @(296, 8),   # "changes,"
    @(308, 10)   # "incentives"
)

foreach ($range in $emphasisRanges) {
    $chars = $tr.Characters($range[0], $range[1])
    $chars.Font.Bold = $true
    $chars.Font.Italic = $true
}
